# Remove leftover boilerplate Class/Table connector+rectangle groups from
# slide 4 ("Tabble and Class boilerplate created"). Each group is one
# Elbow Connector plus the two Rectangles it ties together; the groups
# below are unused duplicates stacked near the top-left corner of the
# slide and are deleted outright.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$namesToDelete = @(
    "Elbow Connector 150",
    "Rectangle 151",
    "Rectangle 152",
    "Elbow Connector 184",
    "Rectangle 185",
    "Rectangle 186",
    "Elbow Connector 187",
    "Rectangle 188",
    "Rectangle 189",
    "Elbow Connector 190",
    "Rectangle 191",
    "Rectangle 192"
)

foreach ($name in $namesToDelete) {
    $s.Shapes.Item($name).Delete()
}
